$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("For plotting")

$data = @{
    2  = @(0.299080303690987, 0.406135218856178, 0.000775248258722531, 0.0273099273380589)
    3  = @(0.267868576908842, 0.371296654292374, 0.00506180723494574, 0.026384713618248)
    4  = @(0.340461036846052, 0.450657693591563, 0.00415226825526129, 0.0281113920269159)
    5  = @(0.0710304522333187, 0.176248277025647, 0.000388496009111305, 0.0268412818347775)
    6  = @(0.0900371495108961, 0.21012414440351, 0.00217168417341075, 0.0306344374726056)
    7  = @(0.117153886279363, 0.239166245291893, 0.00185427280031239, 0.0311256017889107)
    8  = @(0.099555497515911, 0.235134966283164, -0.00290590189549675, 0.0345865991753196)
    9  = @(0.134562464472543, 0.266530509673527, -0.000511604608115842, 0.0336653176533122)
    10 = @(0.152041349939701, 0.280319161842402, -0.00161869347737001, 0.0327239316078319)
    11 = @(0.0924784741064575, 0.369679140761781, 0.000510943633949725, 0.0707144557794192)
    12 = @(0.11277557544517, 0.402457701905151, 0.00285338594647596, 0.0738985016479543)
    13 = @(0.130417595978012, 0.408080237512029, 0.000890941520288434, 0.0708323065137799)
    14 = @(0.0907027098539146, 0.195728202071576, 0.0014250710635687, 0.0267922174024648)
    15 = @(0.0936529159845769, 0.202670604352284, 0.000438344882964851, 0.0278106347876804)
    16 = @(0.108736083040986, 0.219233896406444, 0.000819811316691935, 0.0281882176952698)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Cells.Item($row, 4).Value = $values[0]
    $ws.Cells.Item($row, 5).Value = $values[1]
    $ws.Cells.Item($row, 6).Value = $values[2]
    $ws.Cells.Item($row, 7).Value = $values[3]
}
